$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A268").Value = '0036.024403/2025-01'
$ws.Range("B268").Value = 'medicamentos Soluções orais I'

$ws.Range("A269").Value = '0036.033486/2025-11'
$ws.Range("B269").Value = 'Alimentação Hospitalar - HBAP, HICD, CEMETRON, CRUE, CEREL.'
$ws.Range("C269").Value = 'Licitatório'

$ws.Range("A270").Value = '0036.033230/2025-12'
$ws.Range("B270").Value = 'Fornecimento de alimentação hospitalar - HRC, HEURO'
$ws.Range("C270").Value = 'Licitatório'

$ws.Range("A271").Value = '0036.032746/2025-31'
$ws.Range("B271").Value = 'prestação de serviços de locação de Módulos/Centrais de Compressores de Ar Medicinal - HEPSJP-II.'
$ws.Range("C271").Value = 'Licitatório'

$ws.Range("A272").Value = '0036.032361/2025-74'
$ws.Range("B272").Value = 'Fornecimento de alimentação hospitalar - HRC, HEURO'
$ws.Range("C272").Value = 'Emergencial'

$ws.Range("A273").Value = '0036.020064/2025-86'
$ws.Range("B273").Value = 'GCET - Gerência de Coordenação Estadual de Transplantes. Bolsa Plástica estéril para acondicionamento e isolamento de órgãos'

$ws.Range("A274").Value = '0036.023632/2025-09'
$ws.Range("B274").Value = 'Engenharia Clínica,Serviço de Gerenciamento de Equipamentos Manutenção Corretiva, Preventiva, Preditiva e Calibração dos Equipamentos com Reposição de Peças e Acessórios, Serviço de Gerenciamento de Equipamentos Manutenção Corretiva, Preventiva, Preditiva e Calibração dos Equipamentos com Reposição de Peças e Acessórios - POC'
$ws.Range("C274").Value = 'Licitatório'

$ws.Range("A275").Value = '0036.027089/2025-19'
$ws.Range("B275").Value = 'serviços de manutenção preventiva e corretiva em condicionadores de ar (com fornecimento e reposição de peças, acessórios e componentes eletrônicos) para as unidades administrativa'
$ws.Range("C275").Value = 'Emergencial'

$ws.Range("A276").Value = '0036.017198/2025-10'
$ws.Range("B276").Value = 'Emergencial - serviços de neurologia Cirúrgica e Neurologia Clínica'
$ws.Range("C276").Value = 'Emergencial'

$ws.Range("A277").Value = '0036.016942/2024-88'
$ws.Range("B277").Value = 'Construção de um Espaço de convivência - POC'

[void]$ws.Range("B279").Select()
